$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quotation row appended below the existing data (row 77)
$row = 77

$ws.Cells.Item($row, 1).Value = 45982
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($row, 2).Value = "21,4507"
$ws.Cells.Item($row, 3).Value = "15,7543"
$ws.Cells.Item($row, 4).Value = "15,2104"
$ws.Cells.Item($row, 5).Value = "15,2104"
